{"js": "// The document contains the bold phrase \"DOCX, DOC, PDF, HTML, XPS, R\" followed by\n// a hidden \"_GoBack\" bookmark and then another bold run \"TF and TXT\". The edit\n// removes the bookmark and merges the split text back into a single run reading\n// \"DOCX, DOC, PDF, HTML, XPS, RTF and TXT\".\n\n// 1) Remove the leftover \"_GoBack\" bookmark.\ntry {\n  context.document.deleteBookmark(\"_GoBack\");\n} catch (e) {\n  // Ignore if it doesn't exist (already removed / not present).\n}\nawait context.sync();\n\n// 2) Re-write the (now bookmark-free) phrase so the two adjacent bold runs\n//    collapse into a single run with the full text.\nconst searchResults = context.document.body.search(\n  \"DOCX, DOC, PDF, HTML, XPS, RTF and TXT\",\n  { matchCase: true }\n);\nsearchResults.load(\"text\");\nawait context.sync();\n\nif (searchResults.items.length > 0) {\n  searchResults.items[0].insertText(\n    \"DOCX, DOC, PDF, HTML, XPS, RTF and TXT\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n", "ps1": "# The document contains the bold phrase \"DOCX, DOC, PDF, HTML, XPS, R\" followed\n# by a hidden \"_GoBack\" bookmark and then another bold run \"TF and TXT\". The\n# edit removes the bookmark and merges the split text back into a single run\n# reading \"DOCX, DOC, PDF, HTML, XPS, RTF and TXT\".\n\n$d = $word.ActiveDocument\n\n# 1) Remove the leftover \"_GoBack\" bookmark that splits the bold run in two.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2) Re-write the (now bookmark-free) phrase so the two adjacent bold runs\n#    collapse back into a single run containing the full text.\n$find = $d.Content.Find\n$find.Text = \"DOCX, DOC, PDF, HTML, XPS, RTF and TXT\"\n$find.Replacement.Text = \"DOCX, DOC, PDF, HTML, XPS, RTF and TXT\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n"}
